$d = $word.ActiveDocument

# --- 1) Remove the entire "Meta description" paragraph -------------------
# It's the paragraph made of an empty run + a bold "Meta description" run +
# a normal run with the review blurb. Find it by its text (more robust than
# a hard-coded paragraph index) and delete its whole Range, which removes
# the paragraph (including its trailing paragraph mark) outright.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description:*") {
        $para.Range.Delete()
        break
    }
}

# --- 2) Rework the final paragraph ----------------------------------------
# The last paragraph currently holds the old AI-image prompt ("Create a
# feature image for..."). We need to:
#   a) insert a new bold paragraph right before it:
#        "Play Eureka Reels Blast Superlock Free - Review & Demo"
#   b) replace its own text with the meta-description copy, keeping the
#      paragraph's italic run formatting.
# Doing this as a single OOXML replacement of the last paragraph's Range
# (rather than InsertParagraphBefore + Range.Text=) keeps the new bold
# paragraph from inheriting the old paragraph's italic run formatting, and
# keeps both paragraphs free of any incidental pPr/rsid noise.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$replacement = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Eureka Reels Blast Superlock Free - Review &amp; Demo</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Eureka Reels Blast Superlock, a mining-themed slot game with free spins, Lock It Link feature &amp; Superlock Wheel. Play now for free!</w:t></w:r></w:p>
'@
$lastPara.Range.InsertXML($replacement) | Out-Null
